# release v0.0.2, fixing minor things in printing
#
# - student_info sheet: the image_name cells (G2/G3) both said "test.png";
#   point them at the actual per-student screenshot file names instead,
#   and right-align them like the "Ví dụ" sheet already does.
# - "Ví dụ" sheet: G2/G3 had the file names the other way round relative
#   to student_info - swap them to match.
# - View/printing cleanup: "student_info" becomes the active/selected
#   tab (it was "Ví dụ" before), and the remembered selection on each
#   sheet is refreshed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("student_info")
$ws2 = $wb.Worksheets.Item("Ví dụ")

# --- fix the image file names referenced on each sheet -------------------
$ws1.Range("G2").Value = "quochung.jpg"
$ws1.Range("G3").Value = "vyan.jpg"

$ws2.Range("G2").Value = "vyan.jpg"
$ws2.Range("G3").Value = "quochung.jpg"

# student_info's image cells now get right-aligned, matching "Ví dụ"
$ws1.Range("G2:G3").HorizontalAlignment = -4152

# --- selection / active tab ----------------------------------------------
$ws2.Activate()
$ws2.Range("G2:G3").Select()

$ws1.Activate()
$ws1.Range("C4").Select()
